$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text storage so numeric-looking strings (e.g. "0.999", "10.40")
# are preserved exactly instead of being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.259.39"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "3.119.51"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "593.33"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "157.25"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "3.120.88"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").Value = "5.93"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  -3.87%  "
$ws.Range("D13").Value = "37.23"
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("D15").Value = "3.631.47"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "7.24"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "64.123.98"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "3.114.31"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "481.71"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  -7.68%  "
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "12.97"
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("D26").Value = "81.32"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  +5.14%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("E33").Value = "  -6.41%  "
$ws.Range("D34").Value = "27.48"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").Value = "  -4.61%  "
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  -8.39%  "
$ws.Range("E39").Value = "  -6.25%  "
$ws.Range("D40").Value = "51.14"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "9.23"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "445.13"
$ws.Range("E42").Value = "  -7.24%  "
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("E44").Value = "  -5.23%  "
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "2.843.17"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "130.27"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "25.53"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E51").Value = "  -3.37%  "
